# Updates cryptos list data (prices / 1h volume %, and two swapped-rank
# coin rows) to match the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '94.988.88'
$ws.Range('E2').Value = '  +1.40%  '
# Row 3
$ws.Range('D3').Value = '3.603.89'
$ws.Range('E3').Value = '  +4.53%  '
# Row 4
$ws.Range('E4').Value = '  -0.05%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D5').Value = '235.74'
$ws.Range('E5').Value = '  +0.06%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D6').Value = '658.98'
$ws.Range('E6').Value = '  +5.76%  '
# Row 7
$ws.Range('E7').Value = '  +0.95%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D8').Value = '0.399'
$ws.Range('E8').Value = '  +1.31%  '
# Row 9
$ws.Range('E9').Value = '  -0.06%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D10').Value = '0.988'
$ws.Range('E10').Value = '  -0.98%  '
# Row 11
$ws.Range('D11').Value = '3.598.31'
$ws.Range('E11').Value = '  +4.41%  '
# Row 12
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D12').Value = '0.201'
$ws.Range('E12').Value = '  +0.49%  '
# Row 13
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D13').Value = '42.07'
$ws.Range('E13').Value = '  -3.62%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D14').Value = '6.29'
$ws.Range('E14').Value = '  +0.51%  '
# Row 15
$ws.Range('D15').Value = '4.273.58'
$ws.Range('E15').Value = '  +3.95%  '
# Row 16
$ws.Range('D16').Value = '94.896.47'
$ws.Range('E16').Value = '  +1.44%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D17').Value = '0.0000251'
$ws.Range('E17').Value = '  +0.85%  '
# Row 18
$ws.Range('D18').Value = '3.599.77'
$ws.Range('E18').Value = '  +4.26%  '
# Row 19
$ws.Range('E19').Value = '  -4.49%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D20').Value = '12.77'
$ws.Range('E20').Value = '  +8.76%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D21').Value = '17.86'
$ws.Range('E21').Value = '  -1.28%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D22').Value = '3.53'
$ws.Range('E22').Value = '  +3.75%  '
# Row 23
$ws.Range('E23').Value = '  -2.99%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D24').Value = '502.18'
$ws.Range('E24').Value = '  -1.02%  '
# Row 25
$ws.Range('E25').Value = '  +6.07%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D26').Value = '6.56'
$ws.Range('E26').Value = '  -3.54%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D27').Value = '91.66'
$ws.Range('E27').Value = '  +2.04%  '
# Row 28
$ws.Range('D28').Value = '3.796.73'
$ws.Range('E28').Value = '  +4.35%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D29').Value = '12.43'
$ws.Range('E29').Value = '  +2.24%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D30').Value = '3.11'
$ws.Range('E30').Value = '  +12.87%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  +0.02%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D32').Value = '11.17'
$ws.Range('E32').Value = '  -2.17%  '
# Row 33
$ws.Range('E33').Value = '  -1.81%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D34').Value = '1.01'
$ws.Range('E34').Value = '  +0.36%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D35').Value = '32.03'
$ws.Range('E35').Value = '  +9.05%  '
# Row 36
$ws.Range('E36').Value = '  -1.42%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D37').Value = '0.554'
$ws.Range('E37').Value = '  -0.11%  '
# Row 38
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D38').Value = '565.46'
$ws.Range('E38').Value = '  -1.20%  '
# Row 39
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D39').Value = '8.02'
$ws.Range('E39').Value = '  +6.15%  '
# Row 40
$ws.Range('E40').Value = '  +2.09%  '
# Row 42
$ws.Range('E42').Value = '  +0.08%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D43').Value = '0.909'
$ws.Range('E43').Value = '  -0.17%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D44').Value = '34.87'
$ws.Range('E44').Value = '  +42.03%  '
# Row 45
$ws.Range('E45').Value = '  +1.27%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D46').Value = '23.67'
$ws.Range('E46').Value = '  -0.22%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D47').Value = '5.56'
$ws.Range('E47').Value = '  +0.53%  '
# Row 48
$ws.Range('E48').Value = '  +4.50%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'  # keep numeric-looking text as text
$ws.Range('D49').Value = '0.0410'
$ws.Range('E49').Value = '  -2.32%  '
# Row 50
$ws.Range('E50').Value = '  -1.61%  '
# Row 51
$ws.Range('E51').Value = '  +0.31%  '
